$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.780.92'
$ws.Range('E2').Value = '  +0.59%  '

$ws.Range('D3').Value = '2.305.97'
$ws.Range('E3').Value = '  -0.45%  '

$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.76'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.99%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.42'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.79%  '

$ws.Range('E7').Value = '  -0.75%  '

$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('E9').Value = '  -1.36%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.67'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.06%  '

$ws.Range('E11').Value = '  -0.32%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.50'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.50%  '

$ws.Range('E13').Value = '  +1.47%  '

$ws.Range('E14').Value = '  +2.95%  '

$ws.Range('E15').Value = '  -0.57%  '

$ws.Range('D16').Value = '2.653.37'
$ws.Range('E16').Value = '  -0.55%  '

$ws.Range('D17').Value = '2.300.75'
$ws.Range('E17').Value = '  -0.68%  '

$ws.Range('D18').Value = '42.670.58'
$ws.Range('E18').Value = '  -0.07%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.56'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.28%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.41'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +29.65%  '

$ws.Range('E21').Value = '  -0.05%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.09'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.68%  '

$ws.Range('E23').Value = '  -2.59%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '266.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.69%  '

$ws.Range('E25').Value = '  -2.21%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.34%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.96'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.18%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.28'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.33%  '

$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.81'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +14.04%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.67'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.44%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '37.52'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.42%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '165.68'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.22%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0884'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.26%  '

$ws.Range('E34').Value = '  -3.83%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.58'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.14%  '

$ws.Range('E36').Value = '  -1.81%  '

$ws.Range('E37').Value = '  -1.48%  '

$ws.Range('E38').Value = '  -1.43%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.71'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.64%  '

$ws.Range('E40').Value = '  -2.42%  '

$ws.Range('E41').Value = '  +6.13%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '70.45'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.66%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '96.09'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.06%  '

$ws.Range('E44').Value = '  +0.60%  '

$ws.Range('E45').Value = '  -0.25%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.36'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.70%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '115.62'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.54%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '81.04'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.70%  '

$ws.Range('D49').Value = '1.672.30'
$ws.Range('E49').Value = '  +3.50%  '

$ws.Range('E50').Value = '  -2.11%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.24'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.91%  '
